$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) values
$ws.Range("B2").Value = 289.32168099038631
$ws.Range("C2").Value = 262.52707065191231
$ws.Range("D2").Value = 289.92500927746232
$ws.Range("E2").Value = 257.34769713850034

# Row 3 (STR) values
$ws.Range("B3").Value = 294.73978121257522
$ws.Range("C3").Value = 250.20360560472221
$ws.Range("D3").Value = 308.32675562992677
$ws.Range("E3").Value = 251.96893604029032

# Update the selected range to reflect the new selection noted in the diff
$ws.Range("B1:E3").Select() | Out-Null
